$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.651.18"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.576.61"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.54"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.78"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.04"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.248"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0592"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0891"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.800.71"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.576.48"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.651.24"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.69"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.38"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.16"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  -4.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  +5.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.80"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.03"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.396.87"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0465"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.49"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.962"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.19"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.713.42"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.62"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("E51").Value = "  -0.94%  "
